$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1324.1351
$ws.Range("I15").Value = 1324.1351
$ws.Range("K15").Value = 3972.4053
$ws.Range("M15").Value = -3803.4053
$ws.Range("H19").Value = 2425.7144
$ws.Range("J19").Value = 2746.5
$ws.Range("L19").Value = 2746.5
$ws.Range("N19").Value = -3096.5
$ws.Range("H43").Value = 4045.4546
$ws.Range("I43").Value = 3687.5
$ws.Range("K43").Value = 3687.5
$ws.Range("M43").Value = -3618.5
$ws.Range("H51").Value = 9592.611000000001
$ws.Range("I51").Value = 9205.691999999999
$ws.Range("K51").Value = 9205.691999999999
$ws.Range("M51").Value = -8721.691999999999
$ws.Range("H62").Value = 20836214
$ws.Range("I62").Value = 31251748
$ws.Range("K62").Value = 31251748
$ws.Range("M62").Value = -31251124
$ws.Range("H65").Value = 20836214
$ws.Range("I65").Value = 31251748
$ws.Range("K65").Value = 156258740
$ws.Range("M65").Value = -156255620
$ws.Range("I76").Value = 179065.5
$ws.Range("K76").Value = 179065.5
$ws.Range("M76").Value = -178750.5
$ws.Range("I79").Value = 179065.5
$ws.Range("K79").Value = 179065.5
$ws.Range("M79").Value = -177973.5
$ws.Range("H98").Value = 2250.2693
$ws.Range("J98").Value = 6000
$ws.Range("L98").Value = 6000
$ws.Range("N98").Value = -8996
$ws.Range("H100").Value = 7246.077
$ws.Range("I100").Value = 2490.1428
$ws.Range("J100").Value = 12794.667
$ws.Range("K100").Value = 2490.1428
$ws.Range("L100").Value = 12794.667
$ws.Range("M100").Value = -1949.1428
$ws.Range("N100").Value = -13876.667
$ws.Range("H122").Value = 2250.2693
$ws.Range("J122").Value = 6000
$ws.Range("L122").Value = 18000
$ws.Range("N122").Value = -22900
$ws.Range("H135").Value = 901.11365
$ws.Range("I135").Value = 550.9231
$ws.Range("K135").Value = 4958.3079
$ws.Range("M135").Value = -2423.3079
$ws.Range("H137").Value = 4148.385
$ws.Range("I137").Value = 2041.64
$ws.Range("K137").Value = 6124.92
$ws.Range("M137").Value = -3574.92
$ws.Range("H138").Value = 3358.2827
$ws.Range("J138").Value = 5039.404
$ws.Range("L138").Value = 15118.212
$ws.Range("N138").Value = -25398.212
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2459.254
$ws.Range("I32").Value = 1809.2931
$ws.Range("K32").Value = 1809.2931
$ws.Range("M32").Value = -1522.2931
$ws.Range("H74").Value = 1690.8823
$ws.Range("I74").Value = 1609.3871
$ws.Range("K74").Value = 1609.3871
$ws.Range("M74").Value = -735.3870999999999
$ws.Range("H77").Value = 1690.8823
$ws.Range("I77").Value = 1609.3871
$ws.Range("K77").Value = 8046.9355
$ws.Range("M77").Value = -3678.9355
$ws.Range("H97").Value = 1081.5454
$ws.Range("I97").Value = 1207.2285
$ws.Range("J97").Value = 592.7778
$ws.Range("K97").Value = 1207.2285
$ws.Range("L97").Value = 592.7778
$ws.Range("M97").Value = -711.2284999999999
$ws.Range("N97").Value = -1584.7778
$ws.Range("H122").Value = 4810.8823
$ws.Range("I122").Value = 3621.4285
$ws.Range("K122").Value = 10864.2855
$ws.Range("M122").Value = -8414.2855
$ws.Range("H132").Value = 2785.9092
$ws.Range("I132").Value = 945.5469000000001
$ws.Range("K132").Value = 2836.6407
$ws.Range("M132").Value = -306.6406999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3542.2
$ws.Range("I20").Value = 2566.8
$ws.Range("K20").Value = 2566.8
$ws.Range("M20").Value = -2319.8
$ws.Range("H22").Value = 155
$ws.Range("I22").Value = 155
$ws.Range("K22").Value = 155
$ws.Range("M22").Value = 18
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H94").Value = 883.2632
$ws.Range("I94").Value = 376.8889
$ws.Range("K94").Value = 376.8889
$ws.Range("M94").Value = 74.11110000000002
$ws.Range("H134").Value = 21022.746
$ws.Range("I134").Value = 1916.8889
$ws.Range("K134").Value = 5750.6667
$ws.Range("M134").Value = -3215.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8073.2666
$ws.Range("I16").Value = 2609.9
$ws.Range("J16").Value = 19000
$ws.Range("K16").Value = 2609.9
$ws.Range("L16").Value = 19000
$ws.Range("M16").Value = -2322.9
$ws.Range("N16").Value = -19574
$ws.Range("H88").Value = 12882.5
$ws.Range("J88").Value = 13260
$ws.Range("L88").Value = 13260
$ws.Range("N88").Value = -14072
$ws.Range("H91").Value = 12882.5
$ws.Range("J91").Value = 13260
$ws.Range("L91").Value = 13260
$ws.Range("N91").Value = -16068
$ws.Range("H105").Value = 2693.25
$ws.Range("I105").Value = 2591
$ws.Range("K105").Value = 2591
$ws.Range("M105").Value = -844
$ws.Range("H113").Value = 8073.2666
$ws.Range("I113").Value = 2609.9
$ws.Range("J113").Value = 19000
$ws.Range("K113").Value = 2609.9
$ws.Range("L113").Value = 19000
$ws.Range("M113").Value = -439.9000000000001
$ws.Range("N113").Value = -23340
$ws.Range("H122").Value = 3225.5334
$ws.Range("I122").Value = 2216.6365
$ws.Range("K122").Value = 6649.9095
$ws.Range("M122").Value = -4199.9095
$ws.Range("H132").Value = 3011.3936
$ws.Range("I132").Value = 2205.9111
$ws.Range("J132").Value = 5276.8125
$ws.Range("K132").Value = 6617.7333
$ws.Range("L132").Value = 15830.4375
$ws.Range("M132").Value = -4087.7333
$ws.Range("N132").Value = -20890.4375
$ws.Range("H134").Value = 230838.14
$ws.Range("I134").Value = 2410.7585
$ws.Range("K134").Value = 7232.2755
$ws.Range("M134").Value = -4697.2755
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 36864.684
$ws.Range("I5").Value = 42351.26
$ws.Range("K5").Value = 127053.78
$ws.Range("M5").Value = -126941.78
$ws.Range("H44").Value = 1758.3334
$ws.Range("I44").Value = 83.333336
$ws.Range("J44").Value = 3433.3333
$ws.Range("K44").Value = 250.000008
$ws.Range("L44").Value = 10299.9999
$ws.Range("M44").Value = 147.999992
$ws.Range("N44").Value = -11095.9999
$ws.Range("H92").Value = 801
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 801.2
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 2403.6
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -4899.6
$ws.Range("H112").Value = 5799.5
$ws.Range("I112").Value = 5799.5
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 17398.5
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -16290.5
$ws.Range("N112").ClearContents()
$ws.Range("H135").Value = 36864.684
$ws.Range("I135").Value = 42351.26
$ws.Range("K135").Value = 381161.34
$ws.Range("M135").Value = -378626.34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3203.75
$ws.Range("I22").Value = 2407.625
$ws.Range("J22").Value = 3999.875
$ws.Range("K22").Value = 2407.625
$ws.Range("L22").Value = 3999.875
$ws.Range("M22").Value = -2112.625
$ws.Range("N22").Value = -4589.875
$ws.Range("H27").Value = 3203.75
$ws.Range("I27").Value = 2407.625
$ws.Range("J27").Value = 3999.875
$ws.Range("K27").Value = 2407.625
$ws.Range("L27").Value = 3999.875
$ws.Range("M27").Value = -2300.625
$ws.Range("N27").Value = -4213.875
$ws.Range("H40").Value = 7980
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7980
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7980
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8252
$ws.Range("H132").Value = 4057.6128
$ws.Range("I132").Value = 3453.3076
$ws.Range("K132").Value = 10359.9228
$ws.Range("M132").Value = -7829.9228
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 31210.97
$ws.Range("I132").Value = 1476.3462
$ws.Range("K132").Value = 4429.0386
$ws.Range("M132").Value = -1899.0386
